$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "База данных взята отсюда: <hyperlink>" paragraph together
#    with the empty spacer paragraphs around it, collapsing the original
#    3 paragraphs into a single empty paragraph (keeping sz=28/szCs=28).
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*База данных взята*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # Paragraph right before the "База данных..." paragraph (empty spacer).
    $pBefore = $d.Paragraphs.Item($target - 1)

    # Delete the (empty) paragraph mark before it, merging it away so the
    # "База данных..." text becomes paragraph ($target - 1).
    $rMark1 = $d.Range($pBefore.Range.Start, $pBefore.Range.End)
    $rMark1.Delete()

    # The "База данных..." text is now at index ($target - 1).
    $pText = $d.Paragraphs.Item($target - 1)

    # Delete the paragraph's text (everything except its own paragraph mark).
    $rText = $d.Range($pText.Range.Start, $pText.Range.End - 1)
    $rText.Delete()

    # Delete this (now empty) paragraph's own mark, merging it with the
    # paragraph that follows (the empty spacer that already carries the
    # sz=28/szCs=28 formatting we want to keep).
    $pNowEmpty = $d.Paragraphs.Item($target - 1)
    $rMark2 = $d.Range($pNowEmpty.Range.Start, $pNowEmpty.Range.End)
    $rMark2.Delete()
}

# ---------------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from its old location (after the
#    "...цены на «среднего» держатся на одном уровне." paragraph) onto the
#    newly-collapsed empty paragraph above, mirroring Word's own bookkeeping
#    for the last edit position.
# ---------------------------------------------------------------------------
$oldBookmarkRange = $null
try {
    $oldBm = $d.Bookmarks("_GoBack")
    $oldBookmarkRange = $oldBm.Range.Duplicate
} catch {
    $oldBookmarkRange = $null
}

$newBookmarkHost = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ход решения*") {
        $newBookmarkHost = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($newBookmarkHost -ne $null) {
    $d.Bookmarks.Add("_GoBack", $newBookmarkHost.Range)
}

# ---------------------------------------------------------------------------
# 3) Mark the vacated paragraph (the one that used to hold "_GoBack") as
#    English (US), matching the author's follow-up edit there.
# ---------------------------------------------------------------------------
if ($oldBookmarkRange -ne $null) {
    $oldBookmarkRange.LanguageID = "en-US"
}

# ---------------------------------------------------------------------------
# 4) Drop the two stale <w:lastRenderedPageBreak/> markers that sat in front
#    of "Далее, данные следует..." and "Визуализации" — re-issuing the same
#    text via Find/Replace forces Word to rewrite those runs without the
#    cached page-break marker.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Далее, данные следует проверить на выброс",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Далее, данные следует проверить на выброс", 2) | Out-Null

$d.Content.Find.Execute(
    "Визуализации",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Визуализации", 2) | Out-Null
